$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Cell "name" (User entity, first attribute row) -> split into two runs that
# together read "firstname" ("first" + "name"), keeping the original run's
# rsidRPr/rPr on the first half and giving the new half a plain <w:r> with
# the same rPr, exactly like Word does when new text is typed mid-run.
# ---------------------------------------------------------------------------
$r1 = $d.Content
$found1 = $r1.Find.Execute("name", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not find the 'name' cell to update"
}

$nameXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
    '<w:p w14:paraId="37EAD514" w14:textId="369528A6" w:rsidR="00B83F8B" w:rsidRPr="006106D0" w:rsidRDefault="00B83F8B" w:rsidP="00290C31">' +
      '<w:pPr><w:pStyle w:val="ListParagraph"/><w:spacing w:line="360" w:lineRule="auto"/><w:ind w:left="0"/><w:jc w:val="center"/><w:cnfStyle w:val="000000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:sz w:val="33"/><w:szCs w:val="33"/><w:lang w:val="en-CA"/></w:rPr></w:pPr>' +
      '<w:r w:rsidRPr="006106D0"><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:sz w:val="33"/><w:szCs w:val="33"/><w:lang w:val="en-CA"/></w:rPr><w:t>first</w:t></w:r>' +
      '<w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:sz w:val="33"/><w:szCs w:val="33"/><w:lang w:val="en-CA"/></w:rPr><w:t>name</w:t></w:r>' +
    '</w:p>' +
    '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$r1.InsertXML($nameXml)

# ---------------------------------------------------------------------------
# Cell "surname" (User entity, second attribute row) -> replaced so the text
# reads "lastname", written as two runs: "last" + "name".
# ---------------------------------------------------------------------------
$r2 = $d.Content
$found2 = $r2.Find.Execute("surname", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find the 'surname' cell to update"
}

$surnameXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
    '<w:p w14:paraId="6F77018E" w14:textId="08ABEE12" w:rsidR="00B83F8B" w:rsidRPr="006106D0" w:rsidRDefault="00B83F8B" w:rsidP="00290C31">' +
      '<w:pPr><w:pStyle w:val="ListParagraph"/><w:spacing w:line="360" w:lineRule="auto"/><w:ind w:left="0"/><w:jc w:val="center"/><w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:sz w:val="33"/><w:szCs w:val="33"/><w:lang w:val="en-CA"/></w:rPr></w:pPr>' +
      '<w:r w:rsidRPr="006106D0"><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:sz w:val="33"/><w:szCs w:val="33"/><w:lang w:val="en-CA"/></w:rPr><w:t>last</w:t></w:r>' +
      '<w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:sz w:val="33"/><w:szCs w:val="33"/><w:lang w:val="en-CA"/></w:rPr><w:t>name</w:t></w:r>' +
    '</w:p>' +
    '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$r2.InsertXML($surnameXml)
